$d = $word.ActiveDocument

$d.Content.Find.Execute("56×96=", $true, $false, $false, $false, $false, $true, 1, $false, "78×11=", 2) | Out-Null
$d.Content.Find.Execute("31×86=", $true, $false, $false, $false, $false, $true, 1, $false, "11×36=", 2) | Out-Null
$d.Content.Find.Execute("13×28=", $true, $false, $false, $false, $false, $true, 1, $false, "69×97=", 2) | Out-Null
$d.Content.Find.Execute("94×71=", $true, $false, $false, $false, $false, $true, 1, $false, "23×49=", 2) | Out-Null
$d.Content.Find.Execute("52×50=", $true, $false, $false, $false, $false, $true, 1, $false, "33×34=", 2) | Out-Null
$d.Content.Find.Execute("34×66=", $true, $false, $false, $false, $false, $true, 1, $false, "87×30=", 2) | Out-Null
$d.Content.Find.Execute("67×40=", $true, $false, $false, $false, $false, $true, 1, $false, "18×45=", 2) | Out-Null
$d.Content.Find.Execute("46×90=", $true, $false, $false, $false, $false, $true, 1, $false, "46×40=", 2) | Out-Null
$d.Content.Find.Execute("95×54=", $true, $false, $false, $false, $false, $true, 1, $false, "26×11=", 2) | Out-Null
$d.Content.Find.Execute("67×37=", $true, $false, $false, $false, $false, $true, 1, $false, "68×79=", 2) | Out-Null
$d.Content.Find.Execute("89×98=", $true, $false, $false, $false, $false, $true, 1, $false, "70×80=", 2) | Out-Null
$d.Content.Find.Execute("35×35=", $true, $false, $false, $false, $false, $true, 1, $false, "64×13=", 2) | Out-Null
$d.Content.Find.Execute("99×34=", $true, $false, $false, $false, $false, $true, 1, $false, "34×25=", 2) | Out-Null
$d.Content.Find.Execute("80×80=", $true, $false, $false, $false, $false, $true, 1, $false, "66×17=", 2) | Out-Null
$d.Content.Find.Execute("69×84=", $true, $false, $false, $false, $false, $true, 1, $false, "53×19=", 2) | Out-Null
$d.Content.Find.Execute("55×49=", $true, $false, $false, $false, $false, $true, 1, $false, "36×64=", 2) | Out-Null
$d.Content.Find.Execute("99×59=", $true, $false, $false, $false, $false, $true, 1, $false, "65×53=", 2) | Out-Null
$d.Content.Find.Execute("99×51=", $true, $false, $false, $false, $false, $true, 1, $false, "71×38=", 2) | Out-Null
$d.Content.Find.Execute("45×24=", $true, $false, $false, $false, $false, $true, 1, $false, "84×30=", 2) | Out-Null
$d.Content.Find.Execute("87×40=", $true, $false, $false, $false, $false, $true, 1, $false, "25×46=", 2) | Out-Null
$d.Content.Find.Execute("82×41=", $true, $false, $false, $false, $false, $true, 1, $false, "18×43=", 2) | Out-Null
$d.Content.Find.Execute("40×50=", $true, $false, $false, $false, $false, $true, 1, $false, "32×50=", 2) | Out-Null
$d.Content.Find.Execute("64×14=", $true, $false, $false, $false, $false, $true, 1, $false, "88×89=", 2) | Out-Null
$d.Content.Find.Execute("50×11=", $true, $false, $false, $false, $false, $true, 1, $false, "19×74=", 2) | Out-Null
$d.Content.Find.Execute("41×11=", $true, $false, $false, $false, $false, $true, 1, $false, "19×75=", 2) | Out-Null
